# Fruta / hortaliza, semanal
#
# Inserts one new weekly data row into the "Sandia" sheet at row 181,
# pushing the existing rows 181-205 down to 182-206 (dimension grows from
# A1:R205 to A1:R206), then fills the freshly inserted row with the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 181..205 down to 182..206, leaving a blank row 181 behind.
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new weekly record.
$ws.Cells.Item(181, 1).Value = 9
$ws.Cells.Item(181, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(181, 3).Value = "Metropolitana"
$ws.Cells.Item(181, 4).Value = 44476
$ws.Cells.Item(181, 5).Value = 13
$ws.Cells.Item(181, 6).Value = 100112028
$ws.Cells.Item(181, 7).Value = "Sandia"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 250
$ws.Cells.Item(181, 11).Value = 800
$ws.Cells.Item(181, 12).Value = 1000
$ws.Cells.Item(181, 13).Value = 900
$ws.Cells.Item(181, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(181, 15).Value = "Perú"
$ws.Cells.Item(181, 16).Value = 900
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"
